# Applies the "cryptos list" refresh described by the commit:
#   "Updated cryptos list on Sat Jul 27 21:50:13 UTC 2024 with GitHub Actions"
#
# For each coin row, Price (D) and/or Volume(1h) (E) text is refreshed, and a few
# rows were re-ranked so their Coin/Link/Price/Volume (B/C/D/E) moved to a
# neighboring row (Fetch.AI<->Aptos around rows 32-33, and the
# InjectiveProtocol/OKB/Hedera trio around rows 43-45).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Plain assignment is fine for values Excel will not mistake for a number
    # (e.g. multi-dot thousands-grouped prices, or percentage strings).
    $ws.Range($cell).Value = $text
}

function Set-TextValueForced($cell, $text) {
    # Some Price cells look like plain decimals (e.g. "0.600", "583.21", "0.0692").
    # Assigning those to .Value directly lets Excel coerce them into numbers and
    # drop significant trailing/leading zeros. Writing them as a literal text
    # formula and then pasting the computed value back (values-only) keeps the
    # exact original digit string as a plain string cell, without touching the
    # cell style (unlike toggling NumberFormat to "@").
    $escaped = $text.Replace("""", """""")
    $ws.Range($cell).Formula = "=""" + $escaped + """"
    $ws.Range($cell).Copy() | Out-Null
    $ws.Range($cell).PasteSpecial(-4163) | Out-Null
}

Set-TextValue "D2" "68.803.45"
Set-TextValue "E2" "  +0.99%  "
Set-TextValue "D3" "3.282.42"
Set-TextValue "E3" "  +0.39%  "
Set-TextValue "E4" "  +0.09%  "
Set-TextValueForced "D5" "583.21"
Set-TextValue "E5" "  +0.28%  "
Set-TextValueForced "D6" "185.88"
Set-TextValue "E6" "  +1.79%  "
Set-TextValue "E7" "  +0.08%  "
Set-TextValueForced "D8" "0.600"
Set-TextValue "E8" "  -0.77%  "
Set-TextValue "E9" "  -0.44%  "
Set-TextValue "E10" "  -1.06%  "
Set-TextValueForced "D11" "0.422"
Set-TextValue "E11" "  +1.25%  "
Set-TextValue "D12" "3.858.16"
Set-TextValue "E12" "  +0.64%  "
Set-TextValue "E13" "  -0.06%  "
Set-TextValueForced "D14" "28.46"
Set-TextValue "E14" "  -0.25%  "
Set-TextValue "D15" "68.849.37"
Set-TextValue "E15" "  +1.15%  "
Set-TextValue "E16" "  +1.27%  "
Set-TextValue "D17" "3.272.44"
Set-TextValue "E17" "  +0.08%  "
Set-TextValueForced "D18" "5.86"
Set-TextValue "E18" "  +0.09%  "
Set-TextValueForced "D19" "13.64"
Set-TextValue "E19" "  +0.82%  "
Set-TextValueForced "D20" "395.70"
Set-TextValue "E20" "  +5.18%  "
Set-TextValueForced "D21" "7.72"
Set-TextValue "E21" "  +0.81%  "
Set-TextValue "E22" "  +0.74%  "
Set-TextValue "E23" "  +0.18%  "
Set-TextValueForced "D24" "0.521"
Set-TextValue "E24" "  +1.31%  "
Set-TextValue "E25" "  +0.51%  "
Set-TextValue "E26" "  +4.34%  "
Set-TextValueForced "D27" "9.73"
Set-TextValue "E27" "  +1.00%  "
Set-TextValueForced "D28" "0.999"
Set-TextValue "E28" "  +0.00%  "
Set-TextValue "E29" "  +0.00%  "
Set-TextValueForced "D30" "5.72"
Set-TextValue "E30" "  +0.42%  "
Set-TextValueForced "D31" "23.16"
Set-TextValue "E31" "  +1.34%  "
Set-TextValue "B32" "Aptos"
Set-TextValue "C32" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValueForced "D32" "7.18"
Set-TextValue "E32" "  +3.68%  "
Set-TextValue "B33" "Fetch.AI"
Set-TextValue "C33" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValueForced "D33" "1.30"
Set-TextValue "E33" "  +2.25%  "
Set-TextValueForced "D34" "0.998"
Set-TextValue "E34" "  +0.03%  "
Set-TextValue "E35" "  -0.56%  "
Set-TextValueForced "D36" "163.56"
Set-TextValue "E36" "  +1.08%  "
Set-TextValueForced "D37" "1.97"
Set-TextValue "E37" "  +6.29%  "
Set-TextValueForced "D38" "0.829"
Set-TextValue "E38" "  -2.68%  "
Set-TextValueForced "D39" "26.97"
Set-TextValue "E39" "  +0.34%  "
Set-TextValue "E40" "  -0.46%  "
Set-TextValue "E41" "  -2.41%  "
Set-TextValueForced "D42" "2.56"
Set-TextValue "E42" "  -2.71%  "
Set-TextValue "B43" "OKB"
Set-TextValue "C43" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValueForced "D43" "41.50"
Set-TextValue "E43" "  +1.57%  "
Set-TextValue "B44" "Hedera"
Set-TextValue "C44" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValueForced "D44" "0.0692"
Set-TextValue "E44" "  +1.59%  "
Set-TextValue "B45" "InjectiveProtocol"
Set-TextValue "C45" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValueForced "D45" "25.52"
Set-TextValue "E45" "  -0.91%  "
Set-TextValue "D46" "2.660.16"
Set-TextValue "E46" "  -1.03%  "
Set-TextValueForced "D47" "343.42"
Set-TextValue "E47" "  -2.23%  "
Set-TextValue "E48" "  +0.87%  "
Set-TextValueForced "D49" "32.19"
Set-TextValue "E49" "  +3.00%  "
Set-TextValueForced "D50" "6.37"
Set-TextValue "E50" "  +3.32%  "
Set-TextValueForced "D51" "0.997"
Set-TextValue "E51" "  -0.79%  "

$excel.CutCopyMode = $false

